$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 108, pushing existing rows 108:136 down to 109:137
$ws.Rows.Item(108).Insert()

# Copy the contents/format of the row now at 109 (previously row 108) into the new row 108
$src = $ws.Range("A109:R109")
$dst = $ws.Range("A108:R108")
$src.Copy()
$dst.PasteSpecial()

# Now set the new/changed values for the inserted row 108
$ws.Cells.Item(108, 4).Value = 44663
$ws.Cells.Item(108, 10).Value = 2360
